# "working on data formatting for more clear pipeline"
#
# The abundance values in columns B:D (rows 2-73) are rescaled so that each
# column sums to 100 (i.e. each raw abundance is turned into a percentage of
# its column's total). Column A (protein IDs) and row 1 (header/index row)
# are left untouched.
#
# We do this the "Excel way": compute each column's sum with SUM(), divide
# every cell by its column sum (x100) with a filled-down formula, then copy
# the computed results back onto B:D as static values (PasteSpecial values)
# so the saved file has plain <v> numbers, not formulas - matching the
# target. Helper cells are cleared afterwards so the sheet's used range goes
# back to A1:D73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column totals, stashed in unused helper cells.
$ws.Range("F1").Formula = "=SUM(B2:B73)"
$ws.Range("F2").Formula = "=SUM(C2:C73)"
$ws.Range("F3").Formula = "=SUM(D2:D73)"

# Percent-of-column-total for every row, one helper column per data column.
$ws.Range("G2:G73").Formula = '=B2/$F$1*100'
$ws.Range("H2:H73").Formula = '=C2/$F$2*100'
$ws.Range("I2:I73").Formula = '=D2/$F$3*100'

# Copy the computed percentages back over the original raw values, as
# values only (no formulas left behind).
$ws.Range("G2:I73").Copy()
$ws.Range("B2:D73").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Remove the scratch helper cells used for the computation.
$ws.Range("F1:I73").ClearContents()

# The author's selection moved from A2 to J47 before saving.
$ws.Range("J47").Select() | Out-Null
